$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (the single data row) with the new spot price values
$ws.Range("A2").Value = 46052
$ws.Range("B2").Value = 2.05
$ws.Range("C2").Value = 0.93
$ws.Range("D2").Value = 0.84
$ws.Range("E2").Value = 0.73
$ws.Range("F2").Value = 0.51
$ws.Range("G2").Value = 0.28
$ws.Range("H2").Value = 0.31
$ws.Range("I2").Value = 0.72
$ws.Range("J2").Value = 0.87
$ws.Range("K2").Value = 1.17
$ws.Range("L2").Value = 4.1
$ws.Range("M2").Value = 3.53
$ws.Range("N2").Value = 2.25
$ws.Range("O2").Value = 1.46
$ws.Range("P2").Value = 0.98
$ws.Range("Q2").Value = 1.75
$ws.Range("R2").Value = 2.34
$ws.Range("S2").Value = 1.09
$ws.Range("T2").Value = 1.96
$ws.Range("U2").Value = 18.76
$ws.Range("V2").Value = 21.78
$ws.Range("W2").Value = 20.83
$ws.Range("X2").Value = 15.66
$ws.Range("Y2").Value = 3.62
$ws.Range("Z2").Value = 4.52
$ws.Range("AB2").Value = 15.47
$ws.Range("AD2").Value = 21.3
$ws.Range("AF2").Value = 10.36
$ws.Range("AG2").Value = "0h-23h"
